# ---------------------------------------------------------------------------
# Add "Grand Park via Lake Eleanor" to the Hike Difficulties workbook.
#
# The table on "Hike Difficulties" is sorted alphabetically by Name, so the
# new hike belongs between "Grand Park from Sunrise" (row 14) and
# "Green Lake" (row 15, before the edit) -> the new data lands on row 15
# and every following row shifts down by one.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hike Difficulties")

$tbl = $ws.ListObjects.Item("Table1")

# Row count of the table body before the insert (47 data rows -> last row 48).
$oldLastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1
$newLastRow = $oldLastRow + 1
$insertRow = 15

# Push "Green Lake" ... "White River to Sunrise" down one row and open up
# row 15 for the new hike.
$ws.Rows.Item($insertRow).Insert()

$ws.Cells.Item($insertRow, 1).Value = "Grand Park via Lake Eleanor"
$ws.Cells.Item($insertRow, 2).Value = 9
$ws.Cells.Item($insertRow, 3).Value = 1680
$ws.Cells.Item($insertRow, 4).Value = "moderate"

# Grow the table (and its filter range) to cover the newly inserted row.
$newTableRange = $ws.Range($tbl.Range.Cells.Item(1, 1), $ws.Cells.Item($newLastRow, 4))
$tbl.Resize($newTableRange)

# Re-apply the table's existing alphabetical sort (by Name) across the new range.
$tbl.Sort.SortFields.Clear()
$tbl.Sort.SortFields.Add2($ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($newLastRow, 1)), 0, 1, 0, 2) | Out-Null
$tbl.Sort.Header = 1
$tbl.Sort.Apply()

# Match the saved cursor position from the authored edit.
$ws.Range("D15").Select()
